# Apply targeted numeric value updates to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B3"   = 5.811999999999999
    "A12"  = -21.626
    "B14"  = 6.048
    "B26"  = 6.706
    "A27"  = -21.217
    "B31"  = 6.205
    "A32"  = -21.304
    "B35"  = 7.342000000000001
    "A36"  = -21.141
    "B37"  = 7.848000000000002
    "A38"  = -20.093
    "B45"  = 5.890000000000001
    "A46"  = -21.553
    "B52"  = 5.359
    "A54"  = -21.862
    "A55"  = -22.247
    "A56"  = -22.097
    "B57"  = 5.331999999999999
    "A67"  = -21.588
    "A69"  = -21.721
    "A72"  = -21.567
    "B81"  = 6.809
    "A83"  = -20.146
    "B83"  = 7.326000000000001
    "A86"  = -22.096
    "A91"  = -21.522
    "A93"  = -21.49
    "A99"  = -20.692
    "B100" = 5.517
    "B102" = 7.186
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
